# Apply the "Criado funcionalidade para inserir texto da política de
# privacidade" change:
#   - mark rows 26 and 27 (D26/D27) as "done" (green fill, same as the
#     other completed rows instead of the red "todo" fill)
#   - append a new backlog row (row 36) for the CKEditor / fixed-text
#     forms feature, left marked as "todo" (red fill)
#   - move the visible viewport / active cell down to track the new row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel's Interior.Color is a BGR-packed long (0xBBGGRR).
$red   = 255        # RGB(255,0,0)   -> fill used for "todo"  rows (style s="3")
$green = 5287936     # RGB(0,176,80)  -> fill used for "done"  rows (style s="4")

# --- Mark D26 / D27 as completed (red -> green fill) ------------------
$ws.Range("D26").Interior.Color = $green
$ws.Range("D27").Interior.Color = $green

# --- Append the new backlog row (row 36) -------------------------------
$ws.Range("A36").Value = "Todas"
$ws.Range("B36").Value = "Atualizar para CKEditor formulários de textos fixos no site"
$ws.Range("C36").Value = "Alteração"
$ws.Range("D36").Interior.Color = $red

# --- Update the visible selection / scroll position --------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 25
$win.ScrollColumn = 1
$ws.Range("D28").Select() | Out-Null
